# Auto-generated edit script applying the diff changes to the workbook
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型


# ---- 展览 ----
$ws1.Range("F2").Value = 573
$ws1.Range("F3").Value = 267
$ws1.Range("F4").Value = 606
$ws1.Range("F5").Value = 1428
$ws1.Range("F6").Value = 738
$ws1.Range("F7").Value = 372
$ws1.Range("F8").Value = 61
$ws1.Range("F10").Value = 6668
$ws1.Range("F11").Value = 132
$ws1.Range("F12").Value = 35
$ws1.Range("D14").Value = '松合路2号 钱塘文体中心'
$ws1.Range("F14").Value = 4911
$ws1.Range("I14").Value = '//i1.hdslb.com/bfs/openplatform/202406/czPRn1ve1718875288240.jpeg'
$ws1.Range("F16").Value = 6091
$ws1.Range("F17").Value = 7841
$ws1.Range("F19").Value = 1104
$ws1.Range("F20").Value = 790
$ws1.Range("F21").Value = 4132
$ws1.Range("F22").Value = 606
$ws1.Range("F23").Value = 73
$ws1.Range("F26").Value = 149
$ws1.Range("F27").Value = 1096
$ws1.Range("F28").Value = 26
$ws1.Range("F29").Value = 1545
$ws1.Range("F30").Value = 609
$ws1.Range("F32").Value = 1740
$ws1.Range("F33").Value = 265
$ws1.Range("F34").Value = 2022
$ws1.Range("F35").Value = 253
$ws1.Range("F37").Value = 1306
$ws1.Range("F38").Value = 1329
$ws1.Range("F39").Value = 725
$ws1.Range("F40").Value = 343
$ws1.Range("F41").Value = 3803
$ws1.Range("F42").Value = 167
$ws1.Range("F43").Value = 8
$ws1.Range("F47").Value = 42
$ws1.Range("F48").Value = 118
$ws1.Range("F49").Value = 3993

# ---- 演出 ----
$ws2.Range("F2").Value = 3
$ws2.Range("F9").Value = 20
$ws2.Range("G24").Value = 100
$ws2.Range("F29").Value = 57

# ---- 本地生活 ----
$ws3.Range("F2").Value = 4687

# ---- 全部类型 ----
$ws4.Range("F2").Value = 4687
$ws4.Range("F3").Value = 3
$ws4.Range("F4").Value = 573
$ws4.Range("F8").Value = 267
$ws4.Range("F9").Value = 606
$ws4.Range("F11").Value = 1428
$ws4.Range("F12").Value = 738
$ws4.Range("F13").Value = 62
$ws4.Range("F15").Value = 6668
$ws4.Range("D18").Value = '松合路2号 钱塘文体中心'
$ws4.Range("F18").Value = 4911
$ws4.Range("I18").Value = '//i1.hdslb.com/bfs/openplatform/202406/czPRn1ve1718875288240.jpeg'
$ws4.Range("F19").Value = 6091
$ws4.Range("F20").Value = 6091
$ws4.Range("F21").Value = 7841
$ws4.Range("F23").Value = 1104
$ws4.Range("F24").Value = 790
$ws4.Range("F25").Value = 4132
$ws4.Range("F26").Value = 606
$ws4.Range("F27").Value = 73
$ws4.Range("F29").Value = 149
$ws4.Range("F30").Value = 1096
$ws4.Range("F31").Value = 1545
$ws4.Range("F32").Value = 609
$ws4.Range("F34").Value = 1740
$ws4.Range("F35").Value = 265
$ws4.Range("F36").Value = 2022
$ws4.Range("G39").Value = 100
$ws4.Range("F41").Value = 725
$ws4.Range("F42").Value = 57
$ws4.Range("F43").Value = 343
$ws4.Range("F45").Value = 3803
$ws4.Range("F46").Value = 167
$ws4.Range("F48").Value = 43
$ws4.Range("F49").Value = 118
$ws4.Range("F51").Value = 3993
